$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the date column keeps storing plain text (not auto-converted to a date serial)
$ws.Range("C2:C10").NumberFormat = "@"

# Update region names (rows 2-6 get reordered), date, and values for rows 2-10

$ws.Range("A2").Value = "Mato Grosso"
$ws.Range("C2").Value = "01/10/2024"
$ws.Range("D2").Value = 97.53

$ws.Range("A3").Value = "Santa Catarina"
$ws.Range("C3").Value = "01/10/2024"
$ws.Range("D3").Value = 97.32

$ws.Range("A4").Value = "Rondônia"
$ws.Range("C4").Value = "01/10/2024"
$ws.Range("D4").Value = 97.17

$ws.Range("A5").Value = "Paraná"
$ws.Range("C5").Value = "01/10/2024"
$ws.Range("D5").Value = 96.75

$ws.Range("A6").Value = "Mato Grosso do Sul"
$ws.Range("C6").Value = "01/10/2024"
$ws.Range("D6").Value = 96.27

$ws.Range("C7").Value = "01/10/2024"
$ws.Range("D7").Value = 96.08

$ws.Range("C8").Value = "01/10/2024"
$ws.Range("D8").Value = 91.47

$ws.Range("C9").Value = "01/10/2024"
$ws.Range("D9").Value = 93.83

$ws.Range("C10").Value = "01/10/2024"
$ws.Range("D10").Value = 91.38
